$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 107 (G107=27766)
$ws.Range("H107").Value = 2309.889
$ws.Range("I107").Value = 2309.889
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2309.889
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -389.8890000000001
$ws.Range("N107").ClearContents()
# Row 141 (G141=44161)
$ws.Range("H141").Value = 3761.7727
$ws.Range("I141").Value = 1792.125
$ws.Range("J141").Value = 9014.166999999999
$ws.Range("K141").Value = 5376.375
$ws.Range("L141").Value = 27042.501
$ws.Range("M141").Value = -196.375
$ws.Range("N141").Value = -37402.501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 57 (G57=39767)
$ws.Range("H57").Value = 4000
$ws.Range("I57").Value = 4000
$ws.Range("K57").Value = 4000
$ws.Range("M57").Value = -3516
# Row 61 (G61=43999)
$ws.Range("H61").Value = 6538306
$ws.Range("I61").Value = 9260909
$ws.Range("J61").Value = 4059.2666
$ws.Range("K61").Value = 9260909
$ws.Range("L61").Value = 4059.2666
$ws.Range("M61").Value = -9260697
$ws.Range("N61").Value = -4483.2666
# Row 122 (G122=36168)
$ws.Range("H122").Value = 68358.664
$ws.Range("I122").Value = 112329.89
$ws.Range("J122").Value = 2401.8333
$ws.Range("K122").Value = 336989.67
$ws.Range("L122").Value = 7205.499899999999
$ws.Range("M122").Value = -334539.67
$ws.Range("N122").Value = -12105.4999
# Row 136 (G136=43999)
$ws.Range("H136").Value = 6538306
$ws.Range("I136").Value = 9260909
$ws.Range("J136").Value = 4059.2666
$ws.Range("K136").Value = 27782727
$ws.Range("L136").Value = 12177.7998
$ws.Range("M136").Value = -27780177
$ws.Range("N136").Value = -17277.7998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 99 (G99=19943)
$ws.Range("H99").Value = 1180.5883
$ws.Range("I99").Value = 1050.6923
$ws.Range("J99").Value = 1602.75
$ws.Range("K99").Value = 1050.6923
$ws.Range("L99").Value = 1602.75
$ws.Range("M99").Value = 447.3077000000001
$ws.Range("N99").Value = -4598.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (G16=27691)
$ws.Range("H16").Value = 1802.5625
$ws.Range("I16").Value = 1721.3334
$ws.Range("J16").Value = 1907
$ws.Range("K16").Value = 1721.3334
$ws.Range("L16").Value = 1907
$ws.Range("M16").Value = -1434.3334
$ws.Range("N16").Value = -2481
# Row 68 (G68=10611)
$ws.Range("H68").Value = 23416.166
$ws.Range("J68").Value = 23416.166
$ws.Range("L68").Value = 23416.166
$ws.Range("N68").Value = -24914.166
# Row 71 (G71=10611)
$ws.Range("H71").Value = 23416.166
$ws.Range("J71").Value = 23416.166
$ws.Range("L71").Value = 70248.49800000001
$ws.Range("N71").Value = -77736.49800000001
# Row 82 (G82=10799)
$ws.Range("H82").Value = 82787.336
$ws.Range("J82").Value = 82787.336
$ws.Range("L82").Value = 82787.336
$ws.Range("N82").Value = -83509.336
# Row 85 (G85=10799)
$ws.Range("H85").Value = 82787.336
$ws.Range("J85").Value = 82787.336
$ws.Range("L85").Value = 82787.336
$ws.Range("N85").Value = -85283.336
# Row 99 (G99=36198)
$ws.Range("H99").Value = 2488.4473
$ws.Range("I99").Value = 2369.8
$ws.Range("J99").Value = 2530.8215
$ws.Range("K99").Value = 2369.8
$ws.Range("L99").Value = 2530.8215
$ws.Range("M99").Value = -871.8000000000002
$ws.Range("N99").Value = -5526.8215
# Row 113 (G113=27691)
$ws.Range("H113").Value = 1802.5625
$ws.Range("I113").Value = 1721.3334
$ws.Range("J113").Value = 1907
$ws.Range("K113").Value = 1721.3334
$ws.Range("L113").Value = 1907
$ws.Range("M113").Value = 448.6666
$ws.Range("N113").Value = -6247
# Row 121 (G121=27227)
$ws.Range("H121").Value = 28750
$ws.Range("J121").Value = 28750
$ws.Range("L121").Value = 28750
$ws.Range("N121").Value = -31370
# Row 126 (G126=36198)
$ws.Range("H126").Value = 2488.4473
$ws.Range("I126").Value = 2369.8
$ws.Range("J126").Value = 2530.8215
$ws.Range("K126").Value = 7109.400000000001
$ws.Range("L126").Value = 7592.4645
$ws.Range("M126").Value = -4639.400000000001
$ws.Range("N126").Value = -12532.4645
# Row 133 (G133=43328)
$ws.Range("H133").Value = 34142.6
$ws.Range("I133").Value = 10000
$ws.Range("J133").Value = 50237.668
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 50237.668
$ws.Range("M133").Value = -7470
$ws.Range("N133").Value = -55297.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 100 (G100=19831)
$ws.Range("H100").Value = 10028
$ws.Range("J100").Value = 10028
$ws.Range("L100").Value = 30084
$ws.Range("N100").Value = -31706
# Row 106 (G106=19819)
$ws.Range("H106").Value = 8859
$ws.Range("J106").Value = 8859
$ws.Range("L106").Value = 26577
$ws.Range("N106").Value = -28469
# Row 113 (G113=27843)
$ws.Range("H113").Value = 574.3200000000001
$ws.Range("I113").Value = 617.2727
$ws.Range("K113").Value = 1851.8181
$ws.Range("M113").Value = 318.1819
# Row 134 (G134=44074)
$ws.Range("H134").Value = 6300.49
$ws.Range("I134").Value = 2549.95
$ws.Range("J134").Value = 8720.192999999999
$ws.Range("K134").Value = 7649.849999999999
$ws.Range("L134").Value = 26160.579
$ws.Range("M134").Value = -2579.849999999999
$ws.Range("N134").Value = -36300.579

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 14 (G14=3771)
$ws.Range("H14").Value = 78752.5
$ws.Range("J14").Value = 78752.5
$ws.Range("L14").Value = 78752.5
$ws.Range("N14").Value = -79096.5
# Row 40 (G40=36248)
$ws.Range("H40").Value = 2392.5
$ws.Range("I40").Value = 2471
$ws.Range("K40").Value = 2471
$ws.Range("M40").Value = -2335
# Row 94 (G94=18067)
$ws.Range("H94").Value = 76665
$ws.Range("J94").Value = 76665
$ws.Range("L94").Value = 76665
$ws.Range("N94").Value = -78017
# Row 122 (G122=36247)
$ws.Range("H122").Value = 4913.56
$ws.Range("I122").Value = 2257.1428
$ws.Range("J122").Value = 5946.6113
$ws.Range("K122").Value = 6771.428400000001
$ws.Range("L122").Value = 17839.8339
$ws.Range("M122").Value = -4321.428400000001
$ws.Range("N122").Value = -22739.8339
# Row 131 (G131=35466)
$ws.Range("H131").Value = 10000
$ws.Range("J131").Value = 10000
$ws.Range("L131").Value = 10000
$ws.Range("N131").Value = -20080
# Row 132 (G132=44058)
$ws.Range("H132").Value = 2321.516
$ws.Range("I132").Value = 2030.8
$ws.Range("J132").Value = 3532.8333
$ws.Range("K132").Value = 6092.4
$ws.Range("L132").Value = 10598.4999
$ws.Range("M132").Value = -3562.4
$ws.Range("N132").Value = -15658.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14 (G14=2658)
$ws.Range("H14").Value = 839749.7
$ws.Range("I14").Value = 4750
$ws.Range("J14").Value = 1257249.5
$ws.Range("K14").Value = 4750
$ws.Range("L14").Value = 1257249.5
$ws.Range("M14").Value = -4582
$ws.Range("N14").Value = -1257585.5
# Row 122 (G122=36208)
$ws.Range("H122").Value = 2459.5386
$ws.Range("I122").Value = 2550.4
$ws.Range("K122").Value = 7651.200000000001
$ws.Range("M122").Value = -5201.200000000001
# Row 132 (G132=44029)
$ws.Range("H132").Value = 3878605.8
$ws.Range("I132").Value = 2418
$ws.Range("J132").Value = 15154788
$ws.Range("K132").Value = 7254
$ws.Range("L132").Value = 45464364
$ws.Range("M132").Value = -4724
$ws.Range("N132").Value = -45469424
